# eventsliders.xlsx update
#
# Adds 9 new "Artisan Command" rows to the Commands sheet, right after the
# existing "pidmode(<int>)" row:
#   p-i-d(<p>,<i>,<d>), pidSV(<float>), pidRS(<int>), pidSource(<int>),
#   popup(<msg>[,<int>]), message(<msg>), setCanvasColor(<color>),
#   resetCanvasColor, button(<name>)
# Everything that used to start at row 74 (RC Command section onward through
# the WebSocket Command section) shifts down by 9 rows.

$wb = $excel.ActiveWorkbook

$wsSliders  = $wb.Worksheets.Item("Sliders")
$wsCommands = $wb.Worksheets.Item("Commands")

# ---------------------------------------------------------------------
# 1) Commands sheet: insert 9 blank rows before the old row 74.
# ---------------------------------------------------------------------
$wsCommands.Rows("74:82").Insert()

$newCommands = @(
  @('p-i-d(<p>,<i>,<d>)',      'sets the p-i-d parameters of the PID'),
  @('pidSV(<float>)',          'sets the PID target set value SV'),
  @('pidRS(<int>)',            'activates the PID Ramp-Soak pattern number <n> (1-based!)'),
  @('pidSource(<int>)',        'selects the PID input source with <n> 0: BT, 1: ET (Software PID); <n> in {0,..,3} (Arduino PID)'),
  @('popup(<msg>[,<int>])',    'shows popup with message <msg> which optionally automatically closes after <int> seconds'),
  @('message(<msg>)',          'shows message <msg> in the message line'),
  @('setCanvasColor(<color>)', 'sets canvas color to the RGB-hex <color> like #27f1d3'),
  @('resetCanvasColor',        'resets canvas color'),
  @('button(<name>)',          'activates button <name> from { START, CHARGE, DRY, FCs, FCe, SCs, SCe, DROP, COOL, OFF } ')
)

$r = 74
foreach ($pair in $newCommands) {
    $wsCommands.Cells.Item($r, 2).Value = $pair[0]
    $wsCommands.Cells.Item($r, 3).Value = $pair[1]
    $wsCommands.Rows($r).RowHeight = 13.8
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2) Restore view state: active sheet, selections.
# ---------------------------------------------------------------------
$wsCommands.Activate()
$wsCommands.Range("B74:C82").Select() | Out-Null

$wsSliders.Range("B74:C82").Select() | Out-Null
$wsSliders.Range("B6").Select() | Out-Null

$wsCommands.Activate()
